$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Noticia": fill rows 2-6 (columns C,D,E) with news items + hyperlinks
# ---------------------------------------------------------------------------
$wsNoticia = $wb.Worksheets.Item("Noticia")

# Row 2
$wsNoticia.Range("C2").Value = 'https://federaciondecafeteros.org/wp/listado-noticias/produccion-de-cafe-de-colombia-cae-16-en-agosto/'
$wsNoticia.Range("D2").Value = 'Producción de café de Colombia cae 16% en agosto'
$wsNoticia.Range("E2").Value = @'
La producción de café de Colombia, mayor productor mundial de café arábigo suave lavado, fue de 915.000 sacos de 60 kg de café verde, 16% menos frente a 1,1 millón de sacos registrado en agosto de 2020.
En lo corrido del año (enero-agosto) la producción supera los 7,8 millones de sacos, 8% menos frente a los más de 8,5 millones de sacos del mismo periodo anterior.
En los últimos 12 meses (septiembre 2020-agosto 2021), la producción cayó 7% a casi 13,2 millones de sacos desde los 14,2 millones del mismo periodo anterior.
Y en lo que va del año cafetero (octubre 2020-agosto 2021) la producción registró casi 12,2 millones de sacos, 7% menos comparada con los 13,1 del mismo lapso anterior.
'@
$wsNoticia.Range("E2").WrapText = $true
$wsNoticia.Hyperlinks.Add($wsNoticia.Range("C2"), 'https://federaciondecafeteros.org/wp/listado-noticias/produccion-de-cafe-de-colombia-cae-16-en-agosto/')
$wsNoticia.Rows.Item(2).RowHeight = 360

# Row 3
$wsNoticia.Range("C3").Value = 'https://www.larepublica.co/empresas/las-empresas-de-cafe-han-crecido-hasta-80-en-sus-ventas-en-lo-corrido-de-este-ano-3231889'
$wsNoticia.Range("D3").Value = 'Las empresas de café han crecido hasta 80% en sus ventas en lo corrido de este año'
$wsNoticia.Range("E3").Value = 'Según Euromonitor International, el mercado movió $1,7 billones en 2020, pero crecería hasta $2 billones en ingresos para 2025'
$wsNoticia.Hyperlinks.Add($wsNoticia.Range("C3"), 'https://www.larepublica.co/empresas/las-empresas-de-cafe-han-crecido-hasta-80-en-sus-ventas-en-lo-corrido-de-este-ano-3231889')

# Row 4 (no hyperlink)
$wsNoticia.Range("C4").Value = 'https://www.larepublica.co/indicadores-economicos/commodities/cafe'
$wsNoticia.Range("D4").Value = 'Indicadores de precios en el mercado'
$wsNoticia.Range("E4").Value = 'Indicadores económicos de comodities: Café'

# Row 5
$wsNoticia.Range("C5").Value = 'https://www.larepublica.co/globoeconomia/escasez-de-cafe-de-alta-gama-eleva-los-precios-de-las-variedades-mas-economicas-3225576'
$wsNoticia.Range("D5").Value = 'Escasez de café de alta gama eleva los precios de las variedades más económicas'
$wsNoticia.Range("E5").Value = @'
Los compradores de café están aceptando una escasez mundial de café arábica, el tipo de café de alta gama que prefieren los cafés como Starbucks Corp., y están recurriendo a granos más baratos, lo que hace subir los precios. El café robusta coronó el mayor aumento mensual en más de siete años a medida que la demanda se dispara en medio de múltiples vientos en contra de la oferta. Estos granos se consideran de menor calidad y, a menudo, se utilizan en productos y mezclas de café instantáneo.
'@
$wsNoticia.Hyperlinks.Add($wsNoticia.Range("C5"), 'https://www.larepublica.co/globoeconomia/escasez-de-cafe-de-alta-gama-eleva-los-precios-de-las-variedades-mas-economicas-3225576')

# Row 6 (no hyperlink)
$wsNoticia.Range("C6").Value = 'https://www.eltiempo.com/economia/sectores/analisis-de-ricardo-avila-una-dosis-de-cafeina-para-la-economia-605609'
$wsNoticia.Range("D6").Value = 'Una dosis de cafeína para la economía'
$wsNoticia.Range("E6").Value = @'
A primera vista, el anuncio sobre la muy probable llegada de una nueva masa de aire polar sobre el sur y el centro de Brasil, prevista para esta semana por los servicios meteorológicos, no debería importarles más que a los habitantes de las áreas que pueden ser afectadas por el frío extremo.
'@

# ---------------------------------------------------------------------------
# Sheet "Grupo": fill row 3 (columns C,D,E)
# (order: Nombre, Detalle, Frase -- matches original authoring/shared-string order)
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("Grupo")
$wsGrupo.Range("C3").Value = 'Edgar Orozco Guayara'
$wsGrupo.Range("E3").Value = @'
Técnico en Sistemas y Computaciób - Desarrollador de software durante m´´as de veinte años, astrónomo aficionado, fotógrafo y amante de la ciencia ficción.
'@
$wsGrupo.Range("D3").Value = 'En cada sorbo de café hay una buena idea escondida'
$wsGrupo.Rows.Item(3).RowHeight = 45

# ---------------------------------------------------------------------------
# Sheet "Usuario": fill rows 2-5 (columns B,D) + Activo (E) = 1
# (order matches original authoring/shared-string order: D2,D3,B2,B3,B4,D4,B5,D5)
# ---------------------------------------------------------------------------
$wsUsuario = $wb.Worksheets.Item("Usuario")

$wsUsuario.Range("D2").Value = 'Mónica Angulo'
$wsUsuario.Range("D3").Value = 'Jose Luis Rassa'
$wsUsuario.Range("B2").Value = 'monangu'
$wsUsuario.Range("B3").Value = 'josrass'
$wsUsuario.Range("B4").Value = 'mauruiz'
$wsUsuario.Range("D4").Value = 'Mauricio Ruiz'
$wsUsuario.Range("B5").Value = 'edgoroz'
$wsUsuario.Range("D5").Value = 'Edgar Orozco'

$wsUsuario.Range("E2").Value = 1
$wsUsuario.Range("E3").Value = 1
$wsUsuario.Range("E4").Value = 1
$wsUsuario.Range("E5").Value = 1
